$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 71

$ws.Cells.Item($row, 1).Value = "11/11/2025"
$ws.Cells.Item($row, 2).Value = 0.1975015291293336
$ws.Cells.Item($row, 3).Value = 0.8024984708706664
